$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "1" to "Vani"
$ws.Name = "Vani"

# 2. Replace the ellipsis character "…" with three dots "..." everywhere
#    (this updates the existing shared string in place)
$ws.Cells.Replace("…", "...")

# 3. Clear out the Urban row (row 6) data -> mark as unavailable ("...")
$ws.Range("B6").Value = "..."
$ws.Range("C6").Value = "..."
$ws.Range("D6").Value = "..."
$ws.Range("F6").Value = "..."
$ws.Range("G6").Value = "..."
$ws.Range("H6").Value = "..."

# 4. Clear out most of the Rural row (row 7) data, keep C7 = 5
$ws.Range("B7").Value = "..."
$ws.Range("D7").Value = "..."
$ws.Range("F7").Value = "..."
$ws.Range("G7").Value = "..."
$ws.Range("H7").Value = "..."

# 5. Remove the custom row height on rows 5-7 (reset to default)
$ws.Rows("5:7").AutoFit()

# 6. Delete the blank row 8, shifting the Note row (9) up to row 8
$ws.Rows("8").Delete()
